# Index codes in valuesets for faster lookups
# Adds three new rows of expanded-code data to the "Expansion List" sheet,
# mirroring the existing row (row 13) for formatting, then switches the
# active sheet from "Value Set Info" to "Expansion List".

$wb = $excel.ActiveWorkbook
$wsInfo = $wb.Worksheets.Item("Value Set Info")
$wsExp  = $wb.Worksheets.Item("Expansion List")

# Copy formatting of the existing expanded-code row (13) down onto the
# three new rows (14:16) before filling in values.
$srcRow = $wsExp.Range("A13:F13")
$dstRows = $wsExp.Range("A14:F16")
$srcRow.Copy()
$dstRows.PasteSpecial(-4122) # xlPasteFormats

# Row 14: additional expansion for existing code 10901-7 under a new
# Code System Version (2021-09).
$wsExp.Range("A14").Value = "10901-7"
$wsExp.Range("B14").Value = "Display for 2021-09"
$wsExp.Range("C14").Value = "SNOMEDCT"
$wsExp.Range("D14").Value = "2021-09"
$wsExp.Range("E14").Value = "2.16.840.1.113883.6.96"
$wsExp.Range("F14").Value = "FN"

# Row 15: a new code (10901-8) under SNOMEDCT, same version.
$wsExp.Range("B15").Value = "Display for 10901-8"
$wsExp.Range("A15").Value = "10901-8"
$wsExp.Range("C15").Value = "SNOMEDCT"
$wsExp.Range("D15").Value = "2021-09"
$wsExp.Range("E15").Value = "2.16.840.1.113883.6.96"
$wsExp.Range("F15").Value = "FN"

# Row 16: the same new code (10901-8) but also mapped against LOINC.
$wsExp.Range("A16").Value = "10901-8"
$wsExp.Range("B16").Value = "Display for 10901-8 LOINC"
$wsExp.Range("C16").Value = "LOINC"
$wsExp.Range("D16").Value = "2021-09"
$wsExp.Range("E16").Value = "2.16.840.1.113883.6.1"
$wsExp.Range("F16").Value = "FN"

# Row 16 also carries a slightly shorter custom row height in the source
# workbook.
$wsExp.Rows.Item(16).RowHeight = 14.25

# Fix the selection on "Value Set Info" (no longer the active tab) back to
# its normal cell, and move the frozen-pane selection on "Expansion List"
# down to the newly added last row.
$wsInfo.Range("B3").Select() | Out-Null
$wsExp.Range("A16:XFD16").Select() | Out-Null

# The "Expansion List" sheet becomes the active tab (it was "Value Set
# Info" before).
$wsExp.Activate() | Out-Null
